$d = $word.ActiveDocument

# --- Paragraph: "Jack is almost the opposite of Miles..." ---
$rng = $d.Content
$found = $rng.Find.Execute("Jack is almost the opposite of Miles", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find paragraph 1 anchor text" }
$p = $rng.Paragraphs(1).Range
$xml1 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="68D7801D" w14:textId="77777777" w:rsidR="00D161A1" w:rsidRPr="00670D95" w:rsidRDefault="00D161A1" w:rsidP="00D161A1"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00670D95"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Jack is almost the opposite of Miles; he is </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>really positiv</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>e</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> and optimistic. W</w:t></w:r><w:r w:rsidRPr="00670D95"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">e then meet Stephanie in a winery, serving </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>wine tastings</w:t></w:r><w:r w:rsidRPr="00670D95"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">, and while Jack is being positive, Miles and Stephanie </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>are</w:t></w:r><w:r w:rsidRPr="00670D95"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> negative about the wine. Jack and Stephanie set up a date, almost likes opposites attract each other.</w:t></w:r></w:p>'
$p.InsertXML($xml1)

# --- Paragraph: "There are also other creative elements..." ---
$rng = $d.Content
$found = $rng.Find.Execute("There are also other creative elements", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find paragraph 2 anchor text" }
$p = $rng.Paragraphs(1).Range
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="33D62BDF" w14:textId="77777777" w:rsidR="00D161A1" w:rsidRPr="00670D95" w:rsidRDefault="00D161A1" w:rsidP="00D161A1"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00670D95"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">There are also other creative elements to convey how Miles feels during the movie, if we take a look at the title cards for each day, we can see it’s just a black background with white letters, signifying that Miles is depressed, but the one day where he wakes up with Maya, it’s not a black background, but a colourful view of the local area, signifying Miles isn’t depressed anymore, because of Maya. That day, Maya breaks it off, and the next day </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>we’re</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> back to the black background, signifying Miles is back in a depressed state.</w:t></w:r></w:p>'
$p.InsertXML($xml2)

# --- Paragraph: "In conclusion I really like the red/white contrast..." ---
$rng = $d.Content
$found = $rng.Find.Execute("In conclusion", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find paragraph 3 anchor text" }
$p = $rng.Paragraphs(1).Range
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="3C4912E8" w14:textId="77777777" w:rsidR="00D161A1" w:rsidRDefault="00D161A1" w:rsidP="00D161A1"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r w:rsidRPr="00670D95"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">In conclusion </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">I really like the red/white contrast in the movie, and how like in Miles book, we </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>do not</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> know if he gets together with Maya or what happens to him. </w:t></w:r></w:p>'
$p.InsertXML($xml3)

# --- Paragraph: "I'd recommend watching this movie..." ---
$rng = $d.Content
$found = $rng.Find.Execute("recommend watching this movie", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find paragraph 4 anchor text" }
$p = $rng.Paragraphs(1).Range
$xml4 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w14:paraId="5A671F47" w14:textId="77777777" w:rsidR="00D161A1" w:rsidRPr="00670D95" w:rsidRDefault="00D161A1" w:rsidP="00D161A1"><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>I would</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> recommend watching this movie not only once, but multiple times to really get a sense of how the movie conveys emotions and characters.</w:t></w:r></w:p>'
$p.InsertXML($xml4)

Write-Output "Done"
